$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.425.48'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.675.06'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '644.65'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.96'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.496'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.08'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.450'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.290.17'
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.81'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.664.58'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '69.405.23'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '16.07'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.53'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '467.33'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.92'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.647'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '79.50'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.818.03'
$ws.Range('D24').Style = "Normal"
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  +2.71%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.92'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.11'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  -2.84%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.72'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.02'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.01'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.54%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.96'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.47'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.164'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.92%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.660.96'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '8.44'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.54%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.92'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.37%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '177.57'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.74%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0907'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.926'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '46.64'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.86%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.74'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '27.21'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.37%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.07'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.21%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.86'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.25'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.43%  '
